$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Card (PAN) numbers for data rows 2..10 (column C)
$cards = @(
    "308425000659399",
    "308425000659407",
    "308425000659449",
    "308425000660017",
    "308425000660074",
    "308425000660108",
    "308425000660462",
    "308425000660645",
    "308425000660652"
)

for ($r = 2; $r -le 10; $r++) {
    $idx = $r - 2

    if ($r -ge 4) {
        # brand-new row: clone the formatting of row 2 onto it first
        $ws.Range("A2:K2").Copy() | Out-Null
        $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    }

    # Username / Password - plain shared text
    $ws.Range("A$r").Value = "rambo"
    $ws.Range("B$r").Value = "xyz"

    # CVV (column K) -> shared string, keeps style s=1 (quotePrefix + General).
    # Written before column C so that the shared-string table ends up built in
    # the same order as the original edit ("847" precedes the cardPAN values).
    # Using Formula with a leading apostrophe forces text storage while
    # preserving the cell's existing style (no new style entries created).
    $ws.Range("K$r").Formula = "'847"

    # Card PAN (column C) -> shared string, style s=3 (quotePrefix + Text)
    $ws.Range("C$r").Formula = "'" + $cards[$idx]

    # Top up amount
    $ws.Range("D$r").Value = 10

    $ws.Range("E$r").Value = 4564
    $ws.Range("F$r").Value = 7100

    # Expiry month -> shared string, style s=2 (quotePrefix + numFmt 1)
    $ws.Range("G$r").Formula = "'0000"

    # Expiry year -> shared string, style s=1 (quotePrefix + General)
    $ws.Range("H$r").Formula = "'0004"

    $ws.Range("I$r").Value = "Feb"
    $ws.Range("J$r").Value = 2019
}

# Update the sheet view selection / scroll position to match the target
$ws.Range("B3").Select() | Out-Null
